$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.882.69'
$ws.Cells.Item(2, 5).Value = '  +0.32%  '
$ws.Cells.Item(3, 4).Value = '3.478.91'
$ws.Cells.Item(3, 5).Value = '  +0.99%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '577.07'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.59%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '161.34'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.83%  '
$ws.Cells.Item(7, 5).Value = '  +0.01%  '
$ws.Cells.Item(8, 4).Value = '3.480.89'
$ws.Cells.Item(8, 5).Value = '  +1.00%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.578'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -8.01%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '7.23'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.19%  '
$ws.Cells.Item(11, 5).Value = '  -1.37%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.439'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -1.74%  '
$ws.Cells.Item(13, 4).Value = '4.079.53'
$ws.Cells.Item(13, 5).Value = '  +0.86%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '0.134'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.06%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '27.66'
$c.Style = "Normal"
$ws.Cells.Item(16, 2).Value = 'ShibaInu'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000177'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -7.72%  '
$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(17, 4).Value = '64.950.83'
$ws.Cells.Item(17, 5).Value = '  +0.36%  '
$ws.Cells.Item(18, 4).Value = '3.456.22'
$ws.Cells.Item(18, 5).Value = '  +0.87%  '
$ws.Cells.Item(19, 5).Value = '  -2.95%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '13.88'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.76%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '382.17'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.78%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '7.99'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.71%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '72.86'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.20%  '
$ws.Cells.Item(24, 5).Value = '  +0.20%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '0.534'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -2.78%  '
$ws.Cells.Item(26, 5).Value = '  +2.56%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '9.91'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.10%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '0.179'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +1.27%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.01%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '1.45'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -3.06%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '6.14'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.03%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '2.02'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.54%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '23.42'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.25%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '7.09'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -1.19%  '
$ws.Cells.Item(35, 5).Value = '  -0.52%  '
$ws.Cells.Item(36, 5).Value = '  +0.09%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '1.88'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.45%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.0756'
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -1.90%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '26.98'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +2.14%  '
$ws.Cells.Item(40, 4).Value = '2.879.93'
$ws.Cells.Item(40, 5).Value = '  -2.08%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.817'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +5.75%  '
$ws.Cells.Item(42, 2).Value = 'Filecoin'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '4.54'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.45%  '
$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '6.56'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.61%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '43.01'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.78%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '26.06'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.26%  '
$ws.Cells.Item(46, 5).Value = '  -1.90%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '2.45'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +12.61%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '332.21'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +7.60%  '
$ws.Cells.Item(49, 5).Value = '  -1.22%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.850'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -2.30%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '6.48'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -1.45%  '
